$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing shared string "Elaboration du Gante" -> "Elaboration du Gantt"
$ws.Range("B10").Value = "Elaboration du Gantt"

# Update hours for existing entry on row 10 (0.5 -> 1)
$ws.Range("C10").Value = 1

# Fill in the previously empty row 11
$ws.Range("A11").Value = 43163
$ws.Range("B11").Value = "Elaboration du tableau de répartition des heures"
$ws.Range("C11").Value = 1

# Fill in the previously empty row 12
$ws.Range("A12").Value = 43164
$ws.Range("B12").Value = "Finition du Gantt, Cahier des charges, répartition des heures"
$ws.Range("C12").Value = 1.5
$ws.Rows.Item(12).RowHeight = 30

# Update the active selection shown when the sheet is opened
$ws.Range("E15").Select() | Out-Null
